# Raul's Log - append next day's (9/13/2016) log entries, rows 463-472.
# Minor bug fixes and we are still on locals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Date serial for 9/13/2016 (same value Excel stores internally)
$day = 42626

# Row 463 - Pickup Mic / OSG
$ws.Range("A463").Value = "Pickup Mic"
$ws.Range("B463").Value = $day
$ws.Range("C463").Value = "1600"
$ws.Range("D463").Value = "OSG"
$ws.Range("E463").Value = "1005"
$ws.Range("F463").Value = "Pick up 1 podium mic , 2 desk mics with stands and cables and return to booth behind stage"
$ws.Rows.Item(463).RowHeight = 30

# Row 464 - Demo / OSG
$ws.Range("A464").Value = "Demo"
$ws.Range("B464").Value = $day
$ws.Range("C464").Value = "1630"
$ws.Range("D464").Value = "OSG"
$ws.Range("E464").Value = "2001"

# Row 465 - Demo / ACE
$ws.Range("A465").Value = "Demo"
$ws.Range("B465").Value = $day
$ws.Range("C465").Value = "1900"
$ws.Range("D465").Value = "ACE"
$ws.Range("E465").Value = "010"

# Row 466 - Demo / DB
$ws.Range("A466").Value = "Demo"
$ws.Range("B466").Value = $day
$ws.Range("C466").Value = "1900"
$ws.Range("D466").Value = "DB"
$ws.Range("E466").Value = "0010"

# Row 467 - Demo / DB
$ws.Range("A467").Value = "Demo"
$ws.Range("B467").Value = $day
$ws.Range("C467").Value = "1900"
$ws.Range("D467").Value = "DB"
$ws.Range("E467").Value = "0016"

# Row 468 - Demo / HNE
$ws.Range("A468").Value = "Demo"
$ws.Range("B468").Value = $day
$ws.Range("C468").Value = "1900"
$ws.Range("D468").Value = "HNE"
$ws.Range("E468").Value = "036"

# Row 469 - Demo / HNE
$ws.Range("A469").Value = "Demo"
$ws.Range("B469").Value = $day
$ws.Range("C469").Value = "1900"
$ws.Range("D469").Value = "HNE"
$ws.Range("E469").Value = "401"

# Row 470 - Demo / SSB
$ws.Range("A470").Value = "Demo"
$ws.Range("B470").Value = $day
$ws.Range("C470").Value = "1630"
$ws.Range("D470").Value = "SSB"
$ws.Range("E470").Value = "W141"
$ws.Range("F470").Value = "Using PC, nexk mic  and podium mic (there / test)"

# Row 471 - Operator / SSB
$ws.Range("A471").Value = "Operator"
$ws.Range("B471").Value = $day
$ws.Range("C471").Value = "1715"
$ws.Range("D471").Value = "SSB"
$ws.Range("E471").Value = "W141"
$ws.Range("F471").Value = "Operate event between 17:15-18:00"

# Row 472 - AV Shutdown / SSB
$ws.Range("A472").Value = "AV Shutdown"
$ws.Range("B472").Value = $day
$ws.Range("C472").Value = "2000"
$ws.Range("D472").Value = "SSB"
$ws.Range("E472").Value = "W141"

# Selection state, matching the authored workbook
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 454
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("F473").Select()
